$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-20 Thursday" "2025-11-21 Friday"

Replace-Text "19÷5=" "48÷9="
Replace-Text "19÷2=" "44÷5="
Replace-Text "95÷5=" "98÷5="
Replace-Text "96÷2=" "42÷2="
Replace-Text "68÷7=" "70÷6="

Replace-Text "76÷2=" "34÷3="
Replace-Text "23÷3=" "26÷3="
Replace-Text "23÷2=" "44÷7="
Replace-Text "76÷3=" "53÷8="
Replace-Text "83÷4=" "43÷8="

Replace-Text "48÷3=" "15÷4="
Replace-Text "20÷9=" "79÷4="
Replace-Text "99÷3=" "73÷8="
Replace-Text "35÷7=" "37÷3="
Replace-Text "49÷7=" "34÷8="

Replace-Text "13÷6=" "66÷2="
Replace-Text "19÷9=" "75÷4="
Replace-Text "84÷3=" "65÷3="
Replace-Text "91÷5=" "95÷3="
Replace-Text "66÷7=" "52÷8="

Replace-Text "51÷5=" "60÷5="
Replace-Text "23÷6=" "92÷3="
Replace-Text "47÷3=" "20÷8="
Replace-Text "64÷8=" "61÷5="
Replace-Text "80÷8=" "48÷7="
